$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Periodo Mora" values (swap the two arrears periods)
$ws.Range("E16").Value = "2107"
$ws.Range("E17").Value = "2108"

# Update "Salario Basico" values
$ws.Range("G16").Value = 1000000
$ws.Range("G17").Value = 1000000
